$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Noah"
$ws.Range("C2").Value = "Seligson"
